$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column B (date-as-text) keeps its text representation instead of
# being auto-converted to a date serial number by COM type inference.
$ws.Range("B2:B34").NumberFormat = "@"

$ws.Range("D4").Value = 4.979662106289076
$ws.Range("D7").Value = 2.488389810470949
$ws.Range("D8").Value = 2.628892527610291
$ws.Range("D9").Value = 2.577714676692625
$ws.Range("D10").Value = 2.135650594066437
$ws.Range("D11").Value = 1.914602523047483
$ws.Range("A12").Value = 'Brasil'
$ws.Range("B12").Value = '01/01/2025'
$ws.Range("D12").Value = 0.9151289898246674
$ws.Range("B13").Value = '01/01/2015'
$ws.Range("D13").Value = 1.605626159567575
$ws.Range("B14").Value = '01/01/2016'
$ws.Range("D14").Value = 1.863234877615129
$ws.Range("B15").Value = '01/01/2017'
$ws.Range("D15").Value = 1.99128555278225
$ws.Range("B16").Value = '01/01/2018'
$ws.Range("D16").Value = 1.763794827719464
$ws.Range("B17").Value = '01/01/2019'
$ws.Range("D17").Value = 1.330897691605549
$ws.Range("B18").Value = '01/01/2020'
$ws.Range("D18").Value = 1.060919104018018
$ws.Range("B19").Value = '01/01/2021'
$ws.Range("D19").Value = 0.9992390491236974
$ws.Range("B20").Value = '01/01/2022'
$ws.Range("D20").Value = 1.247186501830941
$ws.Range("B21").Value = '01/01/2023'
$ws.Range("D21").Value = 1.226836841050035
$ws.Range("A22").Value = 'Nordeste'
$ws.Range("B22").Value = '01/01/2024'
$ws.Range("D22").Value = 1.186051307316237
$ws.Range("E22").ClearContents()
$ws.Range("A23").Value = 'Nordeste'
$ws.Range("B23").Value = '01/01/2025'
$ws.Range("D23").Value = 0.6084777557768808
$ws.Range("E23").ClearContents()
$ws.Range("B24").Value = '01/01/2015'
$ws.Range("E24").Value = 23
$ws.Range("B25").Value = '01/01/2016'
$ws.Range("E25").Value = 23
$ws.Range("B26").Value = '01/01/2017'
$ws.Range("E26").Value = 24.5
$ws.Range("B27").Value = '01/01/2018'
$ws.Range("E27").Value = 25
$ws.Range("B28").Value = '01/01/2019'
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 25.5
$ws.Range("B29").Value = '01/01/2020'
$ws.Range("B30").Value = '01/01/2021'
$ws.Range("D30").Value = 0.171036068086038
$ws.Range("E30").Value = 22
$ws.Range("B31").Value = '01/01/2022'
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 24.5
$ws.Range("A32").Value = 'Sergipe'
$ws.Range("B32").Value = '01/01/2023'
$ws.Range("C32").Value = 'Roubo de carga'
$ws.Range("D32").Value = 0.2945573791462633
$ws.Range("E32").Value = 19
$ws.Range("F32").Value = $true
$ws.Range("A33").Value = 'Sergipe'
$ws.Range("B33").Value = '01/01/2024'
$ws.Range("C33").Value = 'Roubo de carga'
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 25
$ws.Range("F33").Value = $true
$ws.Range("A34").Value = 'Sergipe'
$ws.Range("B34").Value = '01/01/2025'
$ws.Range("C34").Value = 'Roubo de carga'
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 25.5
$ws.Range("F34").Value = $true

Write-Output "applied"